$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.857.43"
$ws.Range("E2").Value = "  +7.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.57"
$ws.Range("E3").Value = "  +4.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.54"
$ws.Range("E5").Value = "  +3.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4938"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2787"
$ws.Range("E8").Value = "  +7.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06462"
$ws.Range("E9").Value = "  +4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.808.54"
$ws.Range("E10").Value = "  +4.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.86"
$ws.Range("E11").Value = "  +5.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07111"
$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6494"
$ws.Range("E13").Value = "  +6.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.32"
$ws.Range("E14").Value = "  +9.60%  "

$ws.Range("E15").Value = "  +5.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.832.29"
$ws.Range("E16").Value = "  +8.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007408"
$ws.Range("E18").Value = "  +3.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.26"
$ws.Range("E20").Value = "  +7.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.047.78"
$ws.Range("E21").Value = "  +5.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.602"
$ws.Range("E22").Value = "  +3.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.964"
$ws.Range("E23").Value = "  +4.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.351"
$ws.Range("E24").Value = "  +5.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.85"
$ws.Range("E25").Value = "  +4.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "132.33"
$ws.Range("E26").Value = "  +24.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.60"
$ws.Range("E27").Value = "  +8.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.898"
$ws.Range("E28").Value = "  +6.91%  "

$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.169"
$ws.Range("E30").Value = "  +4.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08369"
$ws.Range("E31").Value = "  +5.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.838"
$ws.Range("E32").Value = "  +3.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04981"

$ws.Range("E34").Value = "  +8.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6797"
$ws.Range("E35").Value = "  +9.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.703"
$ws.Range("E36").Value = "  +4.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.281"
$ws.Range("E37").Value = "  +12.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.762"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9583"
$ws.Range("E39").Value = "  +3.45%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01595"
$ws.Range("E40").Value = "  +6.78%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.044"
$ws.Range("E41").Value = "  +6.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9993"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4101"
$ws.Range("E43").Value = "  +6.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.98"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.232"
$ws.Range("E45").Value = "  +5.50%  "

$ws.Range("E46").Value = "  +5.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05534"
$ws.Range("E47").Value = "  +2.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.132"
$ws.Range("E48").Value = "  +2.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.67"
$ws.Range("E49").Value = "  +5.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3633"
$ws.Range("E50").Value = "  +8.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.309"
$ws.Range("E51").Value = "  +6.08%  "
